$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-key the header row for the new search-key logic.
# Column order: A keyword, B bound-SPU, C bound-name, D type, E related-SPU, F related-SKU, G link
$ws.Range("A1").Value = "关键字"
$ws.Range("C1").Value = "绑定商品名称"
$ws.Range("B1").Value = "绑定商品SPU_ID"
$ws.Range("E1").Value = "关联商品SPU_ID"
$ws.Range("F1").Value = "关联商品SKU_ID"
$ws.Range("G1").Value = "链接"
$ws.Range("D1").Value = "类型(1商品,2链接)"

# Column widths to comfortably fit the new headers (best-fit style).
$ws.Columns("B").ColumnWidth = 15.4
$ws.Columns("C").ColumnWidth = 13.15
$ws.Columns("D").ColumnWidth = 13.15
$ws.Columns("E").ColumnWidth = 15.4
$ws.Columns("F").ColumnWidth = 15.4

$ws.Range("D1").Select() | Out-Null
